$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A, shifting B:F left into A:E
$ws.Range("A:A").Delete()
